$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns keep their original text formatting
# (values such as "565.40" or "1.00" must not be auto-converted to numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.432.66"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.982.20"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.28%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.40"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.01"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.35%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.521"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.972.93"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.38"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +12.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.452"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.67"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.57%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.474.89"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.04"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.981.01"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "59.448.36"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "436.84"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.55"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.718"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.35%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.02"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.30"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.88"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +9.44%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.73"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.36%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +9.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.72"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.19"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0769"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +10.00%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.38%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.34%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.62"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "400.31"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.17%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.737.20"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.250"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.89%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +18.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.33"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.41%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.54%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.27"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.61%  "
